$d = $word.ActiveDocument

$replacements = @(
    @{old="2023-12-25 Monday"; new="2023-12-26 Tuesday"},
    @{old="86÷9=9, 5"; new="49÷4=12, 1"},
    @{old="41÷9=4, 5"; new="29÷5=5, 4"},
    @{old="58÷4=14, 2"; new="25÷6=4, 1"},
    @{old="62÷5=12, 2"; new="22÷8=2, 6"},
    @{old="55÷7=7, 6"; new="92÷5=18, 2"},
    @{old="75÷3=25, 0"; new="97÷4=24, 1"},
    @{old="70÷7=10, 0"; new="98÷8=12, 2"},
    @{old="91÷3=30, 1"; new="87÷6=14, 3"},
    @{old="92÷9=10, 2"; new="58÷5=11, 3"},
    @{old="99÷7=14, 1"; new="71÷3=23, 2"},
    @{old="65÷2=32, 1"; new="80÷6=13, 2"},
    @{old="63÷8=7, 7"; new="66÷2=33, 0"},
    @{old="60÷2=30, 0"; new="77÷4=19, 1"},
    @{old="91÷7=13, 0"; new="24÷8=3, 0"},
    @{old="55÷2=27, 1"; new="26÷7=3, 5"},
    @{old="73÷7=10, 3"; new="36÷4=9, 0"},
    @{old="40÷8=5, 0"; new="93÷5=18, 3"},
    @{old="15÷3=5, 0"; new="25÷4=6, 1"},
    @{old="21÷2=10, 1"; new="89÷4=22, 1"},
    @{old="54÷5=10, 4"; new="37÷9=4, 1"},
    @{old="68÷6=11, 2"; new="58÷8=7, 2"},
    @{old="93÷7=13, 2"; new="45÷2=22, 1"},
    @{old="75÷7=10, 5"; new="59÷5=11, 4"},
    @{old="30÷4=7, 2"; new="72÷2=36, 0"},
    @{old="18÷9=2, 0"; new="50÷4=12, 2"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
